# Update countries & provincias Spain
# Applies the COVID-19 country data refresh:
#   - swap the table rows for Barein/Rumania (Barein now ahead of Rumania)
#   - swap the table rows for Laos/Santa Lucia (Laos now ahead of Santa Lucia)
#   - refresh numeric stats for a handful of countries
#   - bump the "Datos actualizados..." timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (rows keep their numeric data, the label moves) ---
$ws.Range("A48").Value = "Barein"
$ws.Range("A49").Value = "Rumania"

$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 00:11"

# --- Numeric refresh: Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes ---
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 2770080
$ws.Range("C4").Value = 42227
$ws.Range("D4").Value = 1155256
$ws.Range("E4").Value = 1484126
$ws.Range("G4").Value = 576
$ws.Range("H4").Value = 130698

# Row 5 (Brasil)
$ws.Range("B5").Value = 1448753
$ws.Range("C5").Value = 40268
$ws.Range("E5").Value = 598081
$ws.Range("G5").Value = 976
$ws.Range("H5").Value = 60632

# Row 17 (Alemania)
$ws.Range("B17").Value = 196306
$ws.Range("C17").Value = 474
$ws.Range("E17").Value = 7445
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 9061

# Row 23 (Catar)
$ws.Range("B23").Value = 102009
$ws.Range("C23").Value = 4163
$ws.Range("D23").Value = 43407
$ws.Range("E23").Value = 55132
$ws.Range("G23").Value = 136
$ws.Range("H23").Value = 3470

# Row 48 (now Barein)
$ws.Range("B48").Value = 27414
$ws.Range("C48").Value = 656
$ws.Range("D48").Value = 21948
$ws.Range("E48").Value = 5374
$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 92

# Row 49 (now Rumania)
$ws.Range("B49").Value = 27296
$ws.Range("C49").Value = 326
$ws.Range("D49").Value = 19314
$ws.Range("E49").Value = 6315
$ws.Range("G49").Value = 16
$ws.Range("H49").Value = 1667

# Row 72 (Costa de Marfil)
$ws.Range("B72").Value = 8896
$ws.Range("C72").Value = 17
$ws.Range("E72").Value = 507
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 251

# Row 89
$ws.Range("B89").Value = 5154
$ws.Range("C89").Value = 165
$ws.Range("D89").Value = 2722
$ws.Range("E89").Value = 2200
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 232

# Row 100
$ws.Range("D100").Value = 932
$ws.Range("E100").Value = 1902

# Row 171
$ws.Range("B171").Value = 201
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 193
$ws.Range("E171").Value = 7

# Row 179
$ws.Range("D179").Value = 115
$ws.Range("E179").Value = 7

# Row 180
$ws.Range("D180").Value = 89
$ws.Range("E180").Value = 4
